$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.928.89'
$ws.Range('E2').Value = '  -1.68%  '
$ws.Range('D3').Value = '2.484.82'
$ws.Range('E3').Value = '  -1.75%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '585.18'
$ws.Range('E5').Value = '  -1.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '168.01'
$ws.Range('E6').Value = '  -4.91%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.516'
$ws.Range('E8').Value = '  -2.87%  '
$ws.Range('D9').Value = '2.483.81'
$ws.Range('E9').Value = '  -1.76%  '
$ws.Range('E10').Value = '  -4.05%  '
$ws.Range('E11').Value = '  +0.19%  '
$ws.Range('B12').Value = 'Cardano'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.337'
$ws.Range('E12').Value = '  -2.97%  '
$ws.Range('B13').Value = 'Toncoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.94'
$ws.Range('E13').Value = '  -4.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.89'
$ws.Range('E14').Value = '  -3.92%  '
$ws.Range('D15').Value = '2.920.87'
$ws.Range('E15').Value = '  -2.23%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000174'
$ws.Range('E16').Value = '  -2.95%  '
$ws.Range('D17').Value = '66.869.26'
$ws.Range('E17').Value = '  -1.46%  '
$ws.Range('D18').Value = '2.465.30'
$ws.Range('E18').Value = '  -2.70%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.63'
$ws.Range('E19').Value = '  +1.14%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.89'
$ws.Range('E20').Value = '  -1.44%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '362.71'
$ws.Range('E21').Value = '  -0.24%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.05'
$ws.Range('E22').Value = '  -3.96%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.47'
$ws.Range('E23').Value = '  -4.99%  '
$ws.Range('E24').Value = '  +0.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '70.88'
$ws.Range('E25').Value = '  -0.18%  '
$ws.Range('E26').Value = '  -5.79%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.46'
$ws.Range('E27').Value = '  -7.81%  '
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.62%  '
$ws.Range('B29').Value = 'WrappedeETH'
$ws.Range('C29').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D29').Value = '2.603.18'
$ws.Range('E29').Value = '  -2.14%  '
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Value = '0.0₃0932'
$ws.Range('E30').Value = '  -5.85%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.08'
$ws.Range('E31').Value = '  -2.34%  '
$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '515.01'
$ws.Range('E32').Value = '  -5.33%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.83'
$ws.Range('E33').Value = '  -2.43%  '
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.28'
$ws.Range('E34').Value = '  -5.14%  '
$ws.Range('B35').Value = 'FirstDigitalUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.127'
$ws.Range('E36').Value = '  -2.27%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '158.44'
$ws.Range('E37').Value = '  +1.11%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.42'
$ws.Range('E38').Value = '  -3.37%  '
$ws.Range('B39').Value = 'EthereumClassic'
$ws.Range('C39').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.94'
$ws.Range('E39').Value = '  +0.55%  '
$ws.Range('B40').Value = 'WhiteBITCoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.56'
$ws.Range('E40').Value = '  -0.52%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.74'
$ws.Range('E41').Value = '  -3.26%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.95'
$ws.Range('E42').Value = '  -4.88%  '
$ws.Range('B43').Value = 'PolygonEcosystemToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.334'
$ws.Range('E43').Value = '  -6.33%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.45'
$ws.Range('E44').Value = '  -2.66%  '
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '39.24'
$ws.Range('E45').Value = '  -2.02%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '142.68'
$ws.Range('E46').Value = '  -3.08%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.539'
$ws.Range('E47').Value = '  -4.14%  '
$ws.Range('B48').Value = 'Filecoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.60'
$ws.Range('E48').Value = '  -3.64%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0268'
$ws.Range('E49').Value = '  -3.58%  '
$ws.Range('B50').Value = 'Optimism'
$ws.Range('C50').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.65'
$ws.Range('E50').Value = '  -3.19%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0738'
$ws.Range('E51').Value = '  -2.55%  '
